# Update Sage scrape results
# - Row 1 (header): insert a new "The" column at G, shifting the old
#   G1:Y1 topic headers one column to the right (to H1:Z1).
# - Rows 2-11 (data): keep the existing G:Y topic-hit values where they are
#   and append a new Z column (value 0). Replace the scraped paper details
#   (Title/Authors/Year/DOI/Access Type, columns B:F) with the refreshed
#   scrape content.
# - Sheet used range grows from A1:Y11 to A1:Z11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row: shift topic columns G1:Y1 -> H1:Z1, then insert "The" at G1
# ---------------------------------------------------------------------

$oldHeaders = @(
    "Estonia 2007", "Estonia", "Bronze Soldier", "Georgia 2008", "Georgia",
    "Russo-Georgian War", "Stuxnet", "Olympic Games", "Shamoon",
    "Saudi Aramco", "Ras Gas", "Sony", "The Interview", "DNC",
    "Guccifer 2.0", "Democratic National Committee", "Ukrainian Power Grid",
    "BlackEnergy 3", "Sandowrm"
)

# Column Z1 does not exist yet (old sheet only went up to Y1), so give it
# the shared header formatting (bold, bordered, centered) before writing
# the value, reusing the format of an existing header cell so no new style
# entry is created.
$ws.Range("F1").Copy()
$ws.Range("Z1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Walk from the right (Z) back to H so we never overwrite a header before
# it has been copied further right.
for ($i = $oldHeaders.Count - 1; $i -ge 0; $i--) {
    $destCol = 8 + $i   # H=8 .. Z=26
    $ws.Cells.Item(1, $destCol).Value = $oldHeaders[$i]
}

# New header cell - reuse the formatting of an existing header cell (bold,
# bordered, centered) via copy/paste-special so no new style is created.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("G1").Value = "The"

# ---------------------------------------------------------------------
# 2. Data rows 2-11: append the new Z column (topic-hit count, all 0)
# ---------------------------------------------------------------------

foreach ($r in 2..11) {
    $ws.Cells.Item($r, 26).Value = 0   # column Z = 26
}

# ---------------------------------------------------------------------
# 3. Data rows 2-11: refreshed scrape content for Title/Authors/Year/DOI/
#    Access Type (columns B-F)
# ---------------------------------------------------------------------

$rows = @(
    @{ Row=2;  Title="Internet of Things, cybersecurity and governing wicked problems: learning from climate change governance"; Authors="Madeline Carr, Feja Lesniewska"; Year="2020"; DOI="10.1177/0047117820948247"; Access="Open Access" },
    @{ Row=3;  Title="Global versus Local Optimization in Redundancy Resolution of Robotic Manipulators"; Authors="Kazem Kazerounian, Zhaoyu Wang"; Year="1988"; DOI="10.1177/027836498800700501"; Access="Restricted" },
    @{ Row=4;  Title="Cyber and contentious politics: Evidence from the US radical environmental movement"; Authors="Thomas Zeitzoff, Grace Gold"; Year="2024"; DOI="10.1177/00223433231221426"; Access="Restricted" },
    @{ Row=5;  Title="Cyberattacks and public opinion – The effect of uncertainty in guiding preferences"; Authors="Eric Jardine, Nathaniel Porter, Ryan Shandler"; Year="2024"; DOI="10.1177/00223433231218178"; Access="Restricted" },
    @{ Row=6;  Title="A virtual necessity: Some modest steps toward greater cybersecurity"; Authors="Herbert Lin"; Year="2012"; DOI="10.1177/0096340212459039"; Access="Restricted" },
    @{ Row=7;  Title="On 3D simultaneous attack against manoeuvring target with communication delays"; Authors="Zhaohui Liu, Yuezu Lv, Jialing Zhou, Liang Hu"; Year="2020"; DOI="10.1177/1729881419894808"; Access="Open Access" },
    @{ Row=8;  Title="Prioritizing investment in military cyber capability using risk analysis"; Authors="Cayt Rowe, Hossein Seif Zadeh, Ivan L. Garanovich, Li Jiang, Daniel Bilusich, Rick Nunes-Vaz, Anthony Ween"; Year="2019"; DOI="10.1177/1548512917707077"; Access="Restricted" },
    @{ Row=9;  Title="Invisible Digital Front: Can Cyber Attacks Shape Battlefield Events?"; Authors="Nadiya Kostyuk, Yuri M. Zhukov"; Year="2019"; DOI="10.1177/0022002717737138"; Access="Restricted" },
    @{ Row=10; Title="Towards a Chronology of Robotic Art"; Authors="Eduardo Kac"; Year="2001"; DOI="10.1177/135485650100700109"; Access="Restricted" },
    @{ Row=11; Title="A novel ensemble learning approach for fault detection of sensor data in cyber-physical system"; Authors="Ramesh Sneka Nandhini, Ramanathan Lakshmanan"; Year="2023"; DOI="10.3233/JIFS-235809"; Access="Restricted" }
)

foreach ($row in $rows) {
    $r = $row.Row

    $ws.Cells.Item($r, 2).Value = $row.Title    # B - Title
    $ws.Cells.Item($r, 3).Value = $row.Authors  # C - Authors

    # D - Year: force text storage (matches the source file, where this
    # column is plain text) instead of letting Excel auto-convert the
    # numeric-looking string into a number.
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $row.Year
    $ws.Range("A2").Copy()
    $ws.Cells.Item($r, 4).PasteSpecial(-4122)  # xlPasteFormats - restore default style
    $excel.CutCopyMode = $false

    $ws.Cells.Item($r, 5).Value = $row.DOI      # E - DOI
    $ws.Cells.Item($r, 6).Value = $row.Access   # F - Access Type
}
